# Applies the 2024-01-02 cryptos-list refresh (prices / 1h volume%% + two name swaps)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    # Leading apostrophe forces text entry (matches source data, which is
    # all inline/shared strings, even for number-looking prices like "41.00").
    $ws.Range($addr).Value = "'$text"
    $ws.Range($addr).Style = "Normal"
}

# Row 2: Bitcoin
Set-TextCell "D2" "45.137.74"
Set-TextCell "E2" "  +4.89%  "

# Row 3: Ethereum
Set-TextCell "D3" "2.360.51"
Set-TextCell "E3" "  +2.34%  "

# Row 4: TetherUSD
Set-TextCell "E4" "  +0.17%  "

# Row 5: Solana
Set-TextCell "D5" "109.26"
Set-TextCell "E5" "  +3.08%  "

# Row 6: BNB
Set-TextCell "D6" "307.75"
Set-TextCell "E6" "  -1.22%  "

# Row 7: XRP
Set-TextCell "D7" "0.628"
Set-TextCell "E7" "  +0.09%  "

# Row 8: USDC
Set-TextCell "E8" "  -0.23%  "

# Row 9: Cardano
Set-TextCell "D9" "0.613"
Set-TextCell "E9" "  +1.32%  "

# Row 10: Avalanche
Set-TextCell "D10" "41.00"
Set-TextCell "E10" "  +2.20%  "

# Row 11: Dogecoin
Set-TextCell "D11" "0.0913"

# Row 12: Polkadot
Set-TextCell "E12" "  +0.86%  "

# Row 13: TRON
Set-TextCell "E13" "  +1.13%  "

# Row 14: Polygon
Set-TextCell "D14" "0.981"
Set-TextCell "E14" "  -1.06%  "

# Row 15: WrappedliquidstakedEther2.0
Set-TextCell "D15" "2.722.88"
Set-TextCell "E15" "  +2.53%  "

# Row 16: Chainlink
Set-TextCell "D16" "15.41"
Set-TextCell "E16" "  +0.38%  "

# Row 17: WrappedEther
Set-TextCell "D17" "2.348.88"
Set-TextCell "E17" "  +2.06%  "

# Row 18: WrappedBTC
Set-TextCell "D18" "45.109.04"
Set-TextCell "E18" "  +5.18%  "

# Row 19: Uniswap
Set-TextCell "D19" "7.29"
Set-TextCell "E19" "  -2.48%  "

# Row 20: ShibaInu
Set-TextCell "E20" "  +0.93%  "

# Row 21: InternetComputer(DFINITY)
Set-TextCell "D21" "13.06"
Set-TextCell "E21" "  -3.57%  "

# Row 22: Litecoin
Set-TextCell "D22" "73.38"
Set-TextCell "E22" "  -0.22%  "

# Row 23: PancakeSwap
Set-TextCell "D23" "3.44"
Set-TextCell "E23" "  -1.47%  "

# Row 24: BitcoinCash
Set-TextCell "D24" "259.81"
Set-TextCell "E24" "  -3.12%  "

# Row 25: ImmutableX
Set-TextCell "E25" "  +1.10%  "

# Row 26: Dai
Set-TextCell "E26" "  -0.50%  "

# Row 27: Filecoin
Set-TextCell "B27" "Filecoin"
Set-TextCell "C27" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D27" "7.42"
Set-TextCell "E27" "  -6.36%  "

# Row 28: Cosmos
Set-TextCell "B28" "Cosmos"
Set-TextCell "C28" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D28" "11.05"
Set-TextCell "E28" "  +1.04%  "

# Row 29: Toncoin
Set-TextCell "E29" "  +2.61%  "

# Row 30: EthereumClassic
Set-TextCell "D30" "22.37"
Set-TextCell "E30" "  +0.29%  "

# Row 31: Hedera
Set-TextCell "D31" "0.0953"
Set-TextCell "E31" "  +9.65%  "

# Row 32: InjectiveProtocol
Set-TextCell "D32" "37.60"
Set-TextCell "E32" "  -1.56%  "

# Row 33: Monero
Set-TextCell "D33" "169.71"
Set-TextCell "E33" "  +2.42%  "

# Row 34: WEMIXToken
Set-TextCell "D34" "2.91"
Set-TextCell "E34" "  +3.98%  "

# Row 35: Stellar
Set-TextCell "D35" "0.131"
Set-TextCell "E35" "  -0.20%  "

# Row 36: RenderToken
Set-TextCell "B36" "RenderToken"
Set-TextCell "C36" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D36" "4.78"
Set-TextCell "E36" "  +3.05%  "

# Row 37: Kaspa
Set-TextCell "B37" "Kaspa"
Set-TextCell "C37" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D37" "0.115"
Set-TextCell "E37" "  +3.25%  "

# Row 38: LidoDAOToken
Set-TextCell "E38" "  +6.56%  "

# Row 39: NEARProtocol
Set-TextCell "E39" "  +7.22%  "

# Row 40: VeChain
Set-TextCell "D40" "0.0355"
Set-TextCell "E40" "  -0.76%  "

# Row 41: ARBITRUM
Set-TextCell "E41" "  +9.39%  "

# Row 42: BitcoinSV
Set-TextCell "D42" "101.38"
Set-TextCell "E42" "  -4.58%  "

# Row 43: Algorand
Set-TextCell "E43" "  +1.21%  "

# Row 44: Celestia
Set-TextCell "D44" "13.05"
Set-TextCell "E44" "  +6.20%  "

# Row 45: MultiversX
Set-TextCell "D45" "69.65"
Set-TextCell "E45" "  -2.31%  "

# Row 46: FirstDigitalUSD
Set-TextCell "E46" "  -0.48%  "

# Row 47: ordi
Set-TextCell "D47" "81.98"
Set-TextCell "E47" "  +7.79%  "

# Row 48: FraxShare
Set-TextCell "E48" "  +5.60%  "

# Row 49: Aave
Set-TextCell "D49" "112.35"
Set-TextCell "E49" "  +0.93%  "

# Row 50: THORChain
Set-TextCell "E50" "  +6.29%  "

# Row 51: Maker
Set-TextCell "D51" "1.626.05"
Set-TextCell "E51" "  -4.10%  "
